$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.114.40"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "2.233.59"
$ws.Range("E3").Value = "  -4.11%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "2.570.68"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.74%  "
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "2.230.73"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").Value = "41.970.53"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  -4.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.04%  "
$ws.Range("E28").Value = "  +13.41%  "
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0816"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.119"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0297"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("E40").Value = "  -9.40%  "
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.73%  "
$ws.Range("E43").Value = "  -9.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.93%  "
$ws.Range("E45").Value = "  -7.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0988"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  -6.19%  "
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -15.66%  "
$ws.Range("E51").Value = "  -5.32%  "
